$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Sort the data (A2:B28) by column B ascending, keeping the header row fixed.
#    Blank cells in B sort to the bottom, in original relative order.
$dataRange = $ws.Range("A1:B28")
$dataRange.Sort($ws.Range("B1"), 1, $null, $null, 1, 0, 0, 1)

# 2. Add the new column C ("IOOUT" marker column).
#    X for rows whose B id is 1-9, a literal number for id 10, blank afterwards.
$ws.Range("C2").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("C5").Value = "X"
$ws.Range("C6").Value = "X"
$ws.Range("C7").Value = "X"
$ws.Range("C8").Value = "X"
$ws.Range("C9").Value = "X"
$ws.Range("C10").Value = "X"
$ws.Range("C11").Value = 1411

# 3. Header cell C1 + merge with B1, centered.
$ws.Range("C1").Value = ""
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("B1:C1").Merge()

# 4. Column C width.
$ws.Columns.Item(3).ColumnWidth = 2

# 5. AutoFilter across the header row.
$ws.Range("A1:C1").AutoFilter()
$nm = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$1")
$nm.Visible = $false

# 6. Move the active selection to match the authored state.
$ws.Range("D14").Select()

Write-Host "done"
